$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 731.9545000000001
$ws.Range("I39").Value = 297
$ws.Range("J39").Value = 1253.9
$ws.Range("K39").Value = 891
$ws.Range("L39").Value = 3761.7
$ws.Range("M39").Value = -595
$ws.Range("N39").Value = -4353.700000000001
$ws.Range("H40").Value = 9784.857
$ws.Range("I40").Value = 9666
$ws.Range("K40").Value = 9666
$ws.Range("M40").Value = -9491
$ws.Range("H62").Value = 6835.5557
$ws.Range("I62").Value = 6380
$ws.Range("J62").Value = 7200
$ws.Range("K62").Value = 6380
$ws.Range("L62").Value = 7200
$ws.Range("M62").Value = -5756
$ws.Range("N62").Value = -8448
$ws.Range("H65").Value = 6835.5557
$ws.Range("I65").Value = 6380
$ws.Range("J65").Value = 7200
$ws.Range("K65").Value = 31900
$ws.Range("L65").Value = 36000
$ws.Range("M65").Value = -28780
$ws.Range("N65").Value = -42240
$ws.Range("H113").Value = 88001.75
$ws.Range("J113").Value = 6890
$ws.Range("L113").Value = 6890
$ws.Range("N113").Value = -13398
$ws.Range("H118").Value = 807
$ws.Range("I118").Value = 848
$ws.Range("J118").Value = 28
$ws.Range("K118").Value = 2544
$ws.Range("L118").Value = 84
$ws.Range("M118").Value = -887
$ws.Range("N118").Value = -3398
$ws.Range("H125").Value = 6944.5
$ws.Range("I125").Value = 3784.9
$ws.Range("J125").Value = 12210.5
$ws.Range("K125").Value = 34064.1
$ws.Range("L125").Value = 109894.5
$ws.Range("M125").Value = -31604.1
$ws.Range("N125").Value = -114814.5
$ws.Range("H137").Value = 6261.3213
$ws.Range("I137").Value = 7206.619
$ws.Range("K137").Value = 21619.857
$ws.Range("M137").Value = -19069.857
$ws.Range("H138").Value = 4810.25
$ws.Range("I138").Value = 3124.3333
$ws.Range("J138").Value = 5107.7646
$ws.Range("K138").Value = 9372.999899999999
$ws.Range("L138").Value = 15323.2938
$ws.Range("M138").Value = -4232.999899999999
$ws.Range("N138").Value = -25603.2938
$ws.Range("H141").Value = 27997
$ws.Range("I141").Value = 27997
$ws.Range("K141").Value = 83991
$ws.Range("M141").Value = -78811

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16373.732
$ws.Range("J32").Value = 110007
$ws.Range("L32").Value = 110007
$ws.Range("N32").Value = -110581
$ws.Range("H61").Value = 8930.040999999999
$ws.Range("I61").Value = 8038.4634
$ws.Range("K61").Value = 8038.4634
$ws.Range("M61").Value = -7826.4634
$ws.Range("H74").Value = 3036.1738
$ws.Range("I74").Value = 1895.2142
$ws.Range("K74").Value = 1895.2142
$ws.Range("M74").Value = -1021.2142
$ws.Range("H77").Value = 3036.1738
$ws.Range("I77").Value = 1895.2142
$ws.Range("K77").Value = 9476.071
$ws.Range("M77").Value = -5108.071
$ws.Range("H122").Value = 5491.294
$ws.Range("I122").Value = 5326
$ws.Range("K122").Value = 15978
$ws.Range("M122").Value = -13528
$ws.Range("H132").Value = 3521.4348
$ws.Range("I132").Value = 3075.3333
$ws.Range("J132").Value = 4653.846
$ws.Range("K132").Value = 9225.999899999999
$ws.Range("L132").Value = 13961.538
$ws.Range("M132").Value = -6695.999899999999
$ws.Range("N132").Value = -19021.538
$ws.Range("H136").Value = 8930.040999999999
$ws.Range("I136").Value = 8038.4634
$ws.Range("K136").Value = 24115.3902
$ws.Range("M136").Value = -21565.3902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4782.3335
$ws.Range("I31").Value = 3870.6897
$ws.Range("K31").Value = 3870.6897
$ws.Range("M31").Value = -3575.6897
$ws.Range("H34").Value = 4782.3335
$ws.Range("I34").Value = 3870.6897
$ws.Range("K34").Value = 3870.6897
$ws.Range("M34").Value = -3668.6897
$ws.Range("H35").Value = 879.5
$ws.Range("I35").Value = 680.6667
$ws.Range("J35").Value = 998.8
$ws.Range("K35").Value = 680.6667
$ws.Range("L35").Value = 998.8
$ws.Range("M35").Value = -386.6667
$ws.Range("N35").Value = -1586.8
$ws.Range("H58").Value = 3811.0833
$ws.Range("I58").Value = 4633.3335
$ws.Range("K58").Value = 4633.3335
$ws.Range("M58").Value = -4430.3335
$ws.Range("H122").Value = 57951.65
$ws.Range("I122").Value = 81010.57000000001
$ws.Range("J122").Value = 4147.5
$ws.Range("K122").Value = 243031.71
$ws.Range("L122").Value = 12442.5
$ws.Range("M122").Value = -240581.71
$ws.Range("N122").Value = -17342.5
$ws.Range("H136").Value = 3811.0833
$ws.Range("I136").Value = 4633.3335
$ws.Range("K136").Value = 13900.0005
$ws.Range("M136").Value = -11350.0005
$ws.Range("H140").Value = 104319.6
$ws.Range("J140").Value = 104319.6
$ws.Range("L140").Value = 104319.6
$ws.Range("N140").Value = -114679.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 194.33333
$ws.Range("I7").Value = 113.875
$ws.Range("J7").Value = 355.25
$ws.Range("K7").Value = 341.625
$ws.Range("L7").Value = 1065.75
$ws.Range("M7").Value = -229.625
$ws.Range("N7").Value = -1289.75
$ws.Range("H56").Value = 11194.628
$ws.Range("I56").Value = 11194.628
$ws.Range("K56").Value = 11194.628
$ws.Range("M56").Value = -10664.628
$ws.Range("H131").Value = 756089.9
$ws.Range("I131").Value = 1015688.2
$ws.Range("J131").Value = 3254.7
$ws.Range("K131").Value = 3047064.6
$ws.Range("L131").Value = 9764.099999999999
$ws.Range("M131").Value = -3042024.6
$ws.Range("N131").Value = -19844.1
$ws.Range("H138").Value = 1069.2307
$ws.Range("I138").Value = 1069.2307
$ws.Range("K138").Value = 3207.6921
$ws.Range("M138").Value = 1932.3079
$ws.Range("H139").Value = 4853.8184
$ws.Range("I139").Value = 1986.625
$ws.Range("K139").Value = 5959.875
$ws.Range("M139").Value = -819.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7175.175
$ws.Range("J70").Value = 6499.7
$ws.Range("L70").Value = 6499.7
$ws.Range("N70").Value = -7039.7
$ws.Range("H73").Value = 7175.175
$ws.Range("J73").Value = 6499.7
$ws.Range("L73").Value = 6499.7
$ws.Range("N73").Value = -8371.700000000001
$ws.Range("H122").Value = 2098.5833
$ws.Range("I122").Value = 1789.4736
$ws.Range("J122").Value = 3273.2
$ws.Range("K122").Value = 5368.4208
$ws.Range("L122").Value = 9819.599999999999
$ws.Range("M122").Value = -2918.4208
$ws.Range("N122").Value = -14719.6
$ws.Range("H126").Value = 16131.5
$ws.Range("J126").Value = 8124.7144
$ws.Range("L126").Value = 24374.1432
$ws.Range("N126").Value = -29314.1432
$ws.Range("H132").Value = 4817.909
$ws.Range("I132").Value = 3937.375
$ws.Range("J132").Value = 5321.0713
$ws.Range("K132").Value = 11812.125
$ws.Range("L132").Value = 15963.2139
$ws.Range("M132").Value = -9282.125
$ws.Range("N132").Value = -21023.2139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 999
$ws.Range("I55").Value = 999
$ws.Range("J55").Value = 999
$ws.Range("K55").Value = 999
$ws.Range("L55").Value = 999
$ws.Range("M55").Value = -826
$ws.Range("N55").Value = -1345
$ws.Range("H61").Value = 6763.241
$ws.Range("I61").Value = 6001.391
$ws.Range("K61").Value = 6001.391
$ws.Range("M61").Value = -5799.391
$ws.Range("H113").Value = 6763.241
$ws.Range("I113").Value = 6001.391
$ws.Range("K113").Value = 6001.391
$ws.Range("M113").Value = -3831.391
$ws.Range("H122").Value = 7893.25
$ws.Range("I122").Value = 7214.778
$ws.Range("K122").Value = 21644.334
$ws.Range("M122").Value = -19194.334
$ws.Range("H132").Value = 3414095.5
$ws.Range("I132").Value = 58507.4
$ws.Range("K132").Value = 175522.2
$ws.Range("M132").Value = -172992.2
$ws.Range("H136").Value = 11942798
$ws.Range("I136").Value = 15124879
$ws.Range("K136").Value = 45374637
$ws.Range("M136").Value = -45372087
$ws.Range("H138").Value = 149809.75
$ws.Range("J138").Value = 149809.75
$ws.Range("L138").Value = 149809.75
$ws.Range("N138").Value = -160089.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17374.6
$ws.Range("I41").Value = 15989
$ws.Range("J41").Value = 17721
$ws.Range("K41").Value = 15989
$ws.Range("L41").Value = 17721
$ws.Range("M41").Value = -15599
$ws.Range("N41").Value = -18501
$ws.Range("H107").Value = 1287.4762
$ws.Range("J107").Value = 657.6667
$ws.Range("L107").Value = 1973.0001
$ws.Range("N107").Value = -5813.0001
$ws.Range("H122").Value = 5868.8657
$ws.Range("I122").Value = 4666.135
$ws.Range("K122").Value = 13998.405
$ws.Range("M122").Value = -11548.405
$ws.Range("H127").Value = 80000
$ws.Range("J127").Value = 80000
$ws.Range("L127").Value = 80000
$ws.Range("N127").Value = -89920
$ws.Range("H132").Value = 9018.226000000001
$ws.Range("I132").Value = 5820.0835
$ws.Range("K132").Value = 17460.2505
$ws.Range("M132").Value = -14930.2505
